# Highlight the two "to be completed" table-cell answers in yellow.
$d = $word.ActiveDocument
$wdYellow = 7  # WdColorIndex / WdHighlightColorIndex value for yellow

# 1) The "*TODO*" placeholder cell.
$range1 = $d.Content
$found1 = $range1.Find.Execute("*TODO*", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $range1.Font.HighlightColorIndex = $wdYellow
}

# 2) The "Final data frame is exported to ... **Conner to add**" cell.
$range2 = $d.Content
$found2 = $range2.Find.Execute("Final data frame is exported to " + [char]8220 + "final_data.xlsx" + [char]8221 + ". **Conner to add**", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $range2.Font.HighlightColorIndex = $wdYellow
}
